# Swap the order of the names in the "Recorded By" column (G) from
# "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com" for every
# row where that exact value appears.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$lastRow = $ws.UsedRange.Rows.Count
$recordedByRange = $ws.Range("G1:G" + $lastRow)

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$result = $recordedByRange.Replace($oldValue, $newValue)

Write-Host "Replaced '$oldValue' with '$newValue' in range G1:G$lastRow -> success: $result"
